$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterSheet")
$ws.Name = "MasterSheet11"

$ws.Range("E1").Value = 'PS_NUMBER'
$ws.Range("F1").Value = 99003754

$ws.Range("E2").Value = 'Display Name'
$ws.Range("F2").Value = 'Rishab Pankajkumar Ostawal'

$ws.Range("E3").Value = 'Official Email Address'
$ws.Range("F3").Value = 'rishab.ostawal@ltts.com'

$ws.Range("E4").Value = 'Training Hall'
$ws.Range("F4").Value = 'Nalanda'

$ws.Range("E5").Value = 'Floor Number'
$ws.Range("F5").Value = 1

$ws.Range("E6").Value = 'Date Of Joining'
$ws.Range("F6").Value = '8th feb,2021'

$ws.Range("E7").Value = 'Domain'
$ws.Range("F7").Value = 'GT'

$ws.Range("E8").Value = 'Attending Genesis'
$ws.Range("F8").Value = 'Yes'

$ws.Range("E9").Value = 'System Number'
$ws.Range("F9").Value = 15

$ws.Range("E10").Value = 'Team Number'
$ws.Range("F10").Value = 15

$ws.Range("E11").Value = 'BUS NUMBER'
$ws.Range("F11").Value = 1

$ws.Range("E12").Value = 'Working Hours'
$ws.Range("F12").Value = 9

$ws.Range("E13").Value = 'Marks Subject1'
$ws.Range("F13").Value = 54

$ws.Range("E14").Value = 'Marks Subject2'
$ws.Range("F14").Value = 57

$ws.Range("E15").Value = 'Marks Subject3'
$ws.Range("F15").Value = 49

$ws.Range("E16").Value = 'Marks Subject4'
$ws.Range("F16").Value = 37

$ws.Range("E17").Value = 'Marks Subject5'
$ws.Range("F17").Value = 47

$ws.Range("E18").Value = 'Result'
$ws.Range("F18").Value = 'Pass'

$ws.Range("E19").Value = 'Number of clasess attended'
$ws.Range("F19").Value = 40

$ws.Range("E20").Value = 'ATTENDANCE '
$ws.Range("F20").Value = 99

$ws.Range("E21").Value = 'Conduct'
$ws.Range("F21").Value = 'Good'

$ws.Range("E22").Value = 'Experience'
$ws.Range("F22").Value = 'Fresher'

$ws.Range("E23").Value = 'Placement'
$ws.Range("F23").Value = 'College'

$ws.Range("E24").Value = 'Salary '
$ws.Range("F24").Value = 10000

$ws.Range("E25").Value = 'None'

$ws.Range("E26").Value = 'Stream'
$ws.Range("F26").Value = 'ECE'

$ws.Range("E27").Value = 'Address'
$ws.Range("F27").Value = 'Fairfield Marriot'

$ws.Range("E28").Value = 'Area'
$ws.Range("F28").Value = 'Rajajinagar'

$ws.Range("E29").Value = 'Room Number'
$ws.Range("F29").Value = 114

$ws.Range("E30").Value = 'Permanent Address'
$ws.Range("F30").Value = 'AP'

$ws.Range("E31").Value = 'Data1'
$ws.Range("F31").Value = 15

$ws.Range("E32").Value = 'Data2'
$ws.Range("F32").Value = 114
